$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.528436018957346
$ws.Range("C2").Value = 0.435452793834297
$ws.Range("D2").Value = 0.601965601965602
$ws.Range("E2").Value = 0.587301587301587
$ws.Range("F2").Value = 0.452247191011236

$ws.Range("B3").Value = 0.666666666666667
$ws.Range("C3").Value = 0.695364238410596
$ws.Range("D3").Value = 0.760526315789474
$ws.Range("E3").Value = 0.6953125
$ws.Range("F3").Value = 0.569230769230769

$ws.Range("B4").Value = 0.595693779904306
$ws.Range("C4").Value = 0.718095238095238
$ws.Range("D4").Value = 0.786885245901639
$ws.Range("E4").Value = 0.844036697247706
$ws.Range("F4").Value = 0.627764127764128

$ws.Range("B5").Value = 0.792941176470588
$ws.Range("C5").Value = 0.794117647058824
$ws.Range("D5").Value = 0.807228915662651
$ws.Range("E5").Value = 0.816120906801008
$ws.Range("F5").Value = 0.630372492836676

$ws.Range("B6").Value = 0.891304347826087
$ws.Range("C6").Value = 0.847790507364976
$ws.Range("D6").Value = 0.906040268456376
$ws.Range("E6").Value = 0.87378640776699
$ws.Range("F6").Value = 0.6211714132187

$ws.Range("B7").Value = 0.57906976744186
$ws.Range("C7").Value = 0.608247422680412
$ws.Range("D7").Value = 0.618556701030928
$ws.Range("E7").Value = 0.667359667359667
$ws.Range("F7").Value = 0.465809768637532
